# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp string (column A, row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Julio de 2020 a las 14:54"

# --- Re-sorted country labels (data refresh changed the case-count ranking) ---
# Rows 81-83: Bosnia y Herzegovina overtook Bulgaria and Republica de Macedonia
$ws.Cells.Item(81, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(82, 1).Value = "Bulgaria"
$ws.Cells.Item(83, 1).Value = "Republica de Macedonia"

# Rows 210-211: Islas Malvinas / Groenlandia swapped order (tied case counts)
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"
$ws.Cells.Item(211, 1).Value = "Groenlandia"

# --- Updated case counts: row -> (B Casos totales, C Nuevos casos, D Casos activos,
#     E Recuperados, F Casos criticos, G Muertes hoy, H Muertes) ---
$updates = @{
    4   = @(4372056, 217,  2090231, 2131973, 0, 3,  149852)
    16  = @(268934,  1993, 222936,  43238,   0, 27, 2760)
    26  = @(109597,  292,  106328,  3104,    0, 0,  165)
    36  = @(67251,   119,  60492,   6221,    0, 4,  538)
    44  = @(53151,   205,  0,       0,       0, 1,  6141)
    58  = @(30446,   396,  23242,   6781,    0, 6,  423)
    78  = @(13547,   109,  12417,   517,     0, 0,  613)
    80  = @(10621,   152,  3752,    6792,    0, 1,  77)
    81  = @(10498,   183,  4930,    5274,    0, 7,  294)
    82  = @(10427,   0,    5355,    4732,    0, 0,  340)
    83  = @(10213,   127,  5564,    4183,    0, 6,  466)
    98  = @(4881,    24,   3936,    806,     0, 3,  139)
    155 = @(701,     1,    665,     27,      0, 0,  9)
    163 = @(431,     11,   365,     66,      0, 0,  0)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}
